# Update "想去人数" (want-to-go count) figures on the 展览, 演出 and
# 全部类型 sheets to the refreshed values from the latest data pull.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 316
$ws1.Range("F4").Value = 8275
$ws1.Range("F5").Value = 6032
$ws1.Range("F6").Value = 517
$ws1.Range("F7").Value = 102
$ws1.Range("F10").Value = 310
$ws1.Range("F11").Value = 950
$ws1.Range("F12").Value = 81

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 96

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 316
$ws4.Range("F4").Value = 8275
$ws4.Range("F5").Value = 6032
$ws4.Range("F6").Value = 517
$ws4.Range("F7").Value = 102
$ws4.Range("F10").Value = 310
$ws4.Range("F11").Value = 96
$ws4.Range("F15").Value = 950
$ws4.Range("F16").Value = 81
